# Fruta / hortaliza, semanal
# Insert two new rows of weekly price data above the existing last two
# records (old rows 17-18 get pushed down to 19-20, unchanged), and
# populate the freed rows 17-18 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the last two existing data rows (17:18) down to (19:20), shifting
# everything in that range down and leaving rows 17:18 blank (mirrors the
# row-insert semantics seen in the target workbook).
$ws.Range("A17:T18").Insert()

# --- Row 17: new "Primera" quality record for the week of 2021-11-05 ---
$ws.Cells.Item(17, 1).Value = 11
$ws.Cells.Item(17, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(17, 3).Value = "Bíobío"
$ws.Cells.Item(17, 4).Value = 44505
$ws.Cells.Item(17, 5).Value = 8
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100107
$ws.Cells.Item(17, 8).Value = "Otros"
$ws.Cells.Item(17, 9).Value = 100107002
$ws.Cells.Item(17, 10).Value = "Chirimoya"
$ws.Cells.Item(17, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 100
$ws.Cells.Item(17, 14).Value = 2200
$ws.Cells.Item(17, 15).Value = 2200
$ws.Cells.Item(17, 16).Value = 2200
$ws.Cells.Item(17, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(17, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(17, 19).Value = 2200
$ws.Cells.Item(17, 20).Value = 1

# --- Row 18: new "Segunda" quality record for the week of 2021-11-05 ---
$ws.Cells.Item(18, 1).Value = 11
$ws.Cells.Item(18, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(18, 3).Value = "Bíobío"
$ws.Cells.Item(18, 4).Value = 44505
$ws.Cells.Item(18, 5).Value = 8
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100107
$ws.Cells.Item(18, 8).Value = "Otros"
$ws.Cells.Item(18, 9).Value = 100107002
$ws.Cells.Item(18, 10).Value = "Chirimoya"
$ws.Cells.Item(18, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(18, 12).Value = "Segunda"
$ws.Cells.Item(18, 13).Value = 100
$ws.Cells.Item(18, 14).Value = 1800
$ws.Cells.Item(18, 15).Value = 1800
$ws.Cells.Item(18, 16).Value = 1800
$ws.Cells.Item(18, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(18, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 19).Value = 1800
$ws.Cells.Item(18, 20).Value = 1
